# Update vm_pu.xlsx results for Case_2_144 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index mapping: B=2 C=3 D=4 E=5 F=6 I=9 J=10 K=11 L=12 M=13 N=14
$colIndexes = @(2,3,4,5,6,9,10,11,12,13,14)

# Each inner array: row number followed by the new values for columns B,C,D,E,F,I,J,K,L,M,N
$data = @(
    @(2,1.02,1.053484357864784,1.060979717634703,1.063290606883541,1.072084073519835,1.051056521970331,1.058501459241406,1.063705023982717,1.066009645420188,1.074779528711521,1.023228231140827),
    @(3,1.02,1.054536443800366,1.061850561249107,1.064326676754013,1.073132266366866,1.051402066180798,1.05920394849926,1.064390256916568,1.066860142352056,1.07564382721263,1.023468264265792),
    @(4,1.02,1.055217213779193,1.062414068274008,1.064997934734623,1.073811126566452,1.05162434818693,1.059657903271631,1.064833019907203,1.067410712720531,1.076203086955353,1.023623237025454),
    @(5,1.02,1.055503409622099,1.062650969274506,1.06528033499478,1.074096664876954,1.051717482217203,1.059848601373993,1.065019006727405,1.06764223007612,1.076438199616165,1.023688304887652),
    @(6,1.02,1.055551463175962,1.062690746136819,1.065327763156012,1.074144616535498,1.051733101453174,1.059880611946579,1.065050225885619,1.067681106224585,1.076477676048976,1.023699225214687),
    @(7,1.02,1.055221037941303,1.062417233746941,1.065001707382025,1.07381494137217,1.051625593879723,1.059660451956956,1.064835505664838,1.067413806040989,1.076206228543173,1.023624106790508),
    @(8,1.02,1.053839915308801,1.061274020349237,1.063640574971841,1.072438188795726,1.051173571042883,1.058738993468093,1.063936731642789,1.066297024520268,1.075071622067841,1.023309422740934),
    @(9,1.02,1.051406194431245,1.059259651541094,1.061248626545442,1.07001685869461,1.050367035683123,1.057110657301704,1.062348174630427,1.064330981669102,1.073072321059418,1.022752275977614),
    @(10,1.02,1.049783697461858,1.057916843780769,1.059658415795489,1.068405808286787,1.049822619615288,1.056022009615714,1.061285921228101,1.063021561962574,1.071739484591324,1.022379082962389),
    @(11,1.02,1.049081130837304,1.05733542250839,1.058970890027158,1.067708960994055,1.049585287095702,1.05554988121517,1.060825193280956,1.062454874661782,1.07116236187018,1.022217069969079),
    @(12,1.02,1.048820163671613,1.05711946053907,1.058715669757029,1.067450233475624,1.049496891336809,1.055374401002587,1.060653943525799,1.062244427038957,1.070947993575655,1.022156828432455),
    @(13,1.02,1.048876142140792,1.057165784939231,1.058770408243812,1.067505726308455,1.049515863382801,1.055412047113596,1.060690682376896,1.06228956671669,1.070993976264942,1.022169753282199),
    @(14,1.02,1.049059559266278,1.057317570945977,1.058949790227303,1.067687572198698,1.049577985168436,1.055535378224429,1.060811040065839,1.062437478070978,1.071144642103344,1.022212091663319),
    @(15,1.02,1.049172568248357,1.057411091823174,1.059060334254335,1.067799628378757,1.049616228662691,1.055611351922691,1.060885181189874,1.062528617142808,1.071237472406623,1.022238169432145),
    @(16,1.02,1.049830324148056,1.057955431310829,1.059704066676639,1.06845207152574,1.049838336951063,1.056053327698819,1.061316482109467,1.063059177470734,1.071777786427363,1.022389826433508),
    @(17,1.02,1.050242912800786,1.058296887532929,1.060108143243497,1.068861532313258,1.049977232113213,1.056330370163765,1.061586820970441,1.063392064418743,1.072116712377783,1.02248484503281),
    @(18,1.02,1.050483567388717,1.058496055404985,1.060343935172513,1.069100436037379,1.050058093245987,1.056491893259642,1.06174443121659,1.063586260683363,1.072314402552544,1.02254022744116),
    @(19,1.02,1.050565624160909,1.058563966873186,1.060424351233915,1.069181908269734,1.050085638670405,1.056546956404043,1.061798159710847,1.063652481494072,1.072381809857332,1.022559104567107),
    @(20,1.02,1.050198646090053,1.058260252278613,1.060064779213035,1.068817593579059,1.049962345912435,1.056300653482011,1.061557823805811,1.063356345804797,1.072080348788179,1.022474654607468),
    @(21,1.02,1.049005547570649,1.05727287364366,1.058896962356648,1.067634020006755,1.049559698473126,1.055499063342273,1.060775600896906,1.062393920612457,1.071100274746506,1.022199625795316),
    @(22,1.02,1.048255383697222,1.056652091505477,1.058163620420923,1.066890511850917,1.049305150097686,1.054994431307603,1.060283121679829,1.06178906772364,1.070484067883747,1.022026341565875),
    @(23,1.02,1.048653061369095,1.056981177686813,1.058552292426874,1.067284597771107,1.049440222567919,1.055262007012668,1.060544257201856,1.062109686782972,1.070810730318899,1.022118237162184),
    @(24,1.02,1.050218648315282,1.058276806166113,1.060084373239966,1.068837447378619,1.049969072821757,1.056314081390887,1.061570926604388,1.063372485419452,1.07209677992773,1.022479259346638),
    @(25,1.02,1.052035371622193,1.059780397492411,1.061866224977973,1.070642273396276,1.050576730366728,1.057532166518545,1.062759422415975,1.064839027449322,1.073589184103798,1.022896622938767),
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $colIndexes.Length; $i++) {
        $c = $colIndexes[$i]
        $value = $rowData[$i + 1]
        $ws.Cells.Item($r, $c).Value = $value
    }
}

